$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 1).Value = "Love_Is"
$ws.Cells.Item(10, 2).Value = "Love is"
$ws.Cells.Item(10, 3).Value = "Kärlek är"
$ws.Cells.Item(10, 4).Value = "Need review"
$ws.Cells.Item(10, 5).Font.Bold = $false
